# chore: update Sheets via scheduled runner
# Refresh cached market-board figures (currentAveragePrice / LevePrice / LeveProfit
# columns H:N) for the rows whose underlying market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3143
$ws.Range("I32").Value = 2924.25
$ws.Range("J32").Value = 3288.8333
$ws.Range("K32").Value = 2924.25
$ws.Range("L32").Value = 3288.8333
$ws.Range("M32").Value = -2598.25
$ws.Range("N32").Value = -3940.8333

$ws.Range("H53").Value = 282.72223
$ws.Range("I53").Value = 193.46153
$ws.Range("K53").Value = 193.46153
$ws.Range("M53").Value = 443.53847

$ws.Range("H86").Value = 13187.857
$ws.Range("J86").Value = 14993.75
$ws.Range("L86").Value = 14993.75
$ws.Range("N86").Value = -17239.75

$ws.Range("H89").Value = 13187.857
$ws.Range("J89").Value = 14993.75
$ws.Range("L89").Value = 74968.75
$ws.Range("N89").Value = -86200.75

$ws.Range("H116").Value = 3635.1304
$ws.Range("J116").Value = 4354.5454
$ws.Range("L116").Value = 4354.5454
$ws.Range("N116").Value = -11238.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1500
$ws.Range("I74").Value = 1000
$ws.Range("K74").Value = 1000
$ws.Range("M74").Value = -126

$ws.Range("H77").Value = 1500
$ws.Range("I77").Value = 1000
$ws.Range("K77").Value = 5000
$ws.Range("M77").Value = -632

$ws.Range("H97").Value = 652.8570999999999
$ws.Range("I97").Value = 661
$ws.Range("K97").Value = 661
$ws.Range("M97").Value = -165

$ws.Range("H122").Value = 28196.916
$ws.Range("I122").Value = 30669.273
$ws.Range("J122").Value = 1001
$ws.Range("K122").Value = 92007.819
$ws.Range("L122").Value = 3003
$ws.Range("M122").Value = -89557.819
$ws.Range("N122").Value = -7903

$ws.Range("H132").Value = 2099
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 149
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H99").Value = 3134.8333
$ws.Range("I99").Value = 3061.8
$ws.Range("K99").Value = 3061.8
$ws.Range("M99").Value = -1563.8

$ws.Range("H134").Value = 3509
$ws.Range("I134").Value = 3676.6667
$ws.Range("K134").Value = 11030.0001
$ws.Range("M134").Value = -8495.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1581.5
$ws.Range("I31").Value = 1486.4
$ws.Range("J31").Value = 1602.174
$ws.Range("K31").Value = 1486.4
$ws.Range("L31").Value = 1602.174
$ws.Range("M31").Value = -1191.4
$ws.Range("N31").Value = -2192.174

$ws.Range("H34").Value = 1581.5
$ws.Range("I34").Value = 1486.4
$ws.Range("J34").Value = 1602.174
$ws.Range("K34").Value = 1486.4
$ws.Range("L34").Value = 1602.174
$ws.Range("M34").Value = -1284.4
$ws.Range("N34").Value = -2006.174

$ws.Range("H42").Value = 46333.332
$ws.Range("J42").Value = 45000
$ws.Range("L42").Value = 45000
$ws.Range("N42").Value = -46186

$ws.Range("H55").Value = 57214.2
$ws.Range("I55").Value = 68691
$ws.Range("K55").Value = 68691
$ws.Range("M55").Value = -68376

$ws.Range("H58").Value = 2473.9333
$ws.Range("I58").Value = 2044.5834
$ws.Range("K58").Value = 2044.5834
$ws.Range("M58").Value = -1841.5834

$ws.Range("H62").Value = 2599
$ws.Range("J62").Value = 2599
$ws.Range("L62").Value = 2599
$ws.Range("N62").Value = -3847

$ws.Range("H65").Value = 2599
$ws.Range("J65").Value = 2599
$ws.Range("L65").Value = 12995
$ws.Range("N65").Value = -19235

$ws.Range("H107").Value = 895.75
$ws.Range("I107").Value = 857.8570999999999
$ws.Range("K107").Value = 857.8570999999999
$ws.Range("M107").Value = 1062.1429

$ws.Range("H134").Value = 1665
$ws.Range("I134").Value = 1635.8948
$ws.Range("K134").Value = 4907.6844
$ws.Range("M134").Value = -2372.6844

$ws.Range("H136").Value = 2473.9333
$ws.Range("I136").Value = 2044.5834
$ws.Range("K136").Value = 6133.7502
$ws.Range("M136").Value = -3583.7502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 133.55
$ws.Range("I2").Value = 134.82353
$ws.Range("K2").Value = 134.82353
$ws.Range("M2").Value = -21.82353000000001

$ws.Range("H59").Value = 6000
$ws.Range("I59").Value = 6000
$ws.Range("J59").Value = 6000
$ws.Range("K59").Value = 6000
$ws.Range("L59").Value = 6000
$ws.Range("M59").Value = -5417
$ws.Range("N59").Value = -7166

$ws.Range("H122").Value = 2583.0833
$ws.Range("I122").Value = 2425
$ws.Range("K122").Value = 7275
$ws.Range("M122").Value = -4825

$ws.Range("H132").Value = 4196.9
$ws.Range("I132").Value = 4329.8887
$ws.Range("K132").Value = 12989.6661
$ws.Range("M132").Value = -10459.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 24995
$ws.Range("J59").Value = 24995
$ws.Range("L59").Value = 24995
$ws.Range("N59").Value = -26303

$ws.Range("H132").Value = 3160.9167
$ws.Range("I132").Value = 2875.4285
$ws.Range("J132").Value = 3560.6
$ws.Range("K132").Value = 8626.2855
$ws.Range("L132").Value = 10681.8
$ws.Range("M132").Value = -6096.2855
$ws.Range("N132").Value = -15741.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12496.75
$ws.Range("I62").Value = 14995.667
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 14995.667
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -14371.667
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 12496.75
$ws.Range("I65").Value = 14995.667
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 74978.33499999999
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -71858.33499999999
$ws.Range("N65").Value = -31240

$ws.Range("H75").Value = 24559
$ws.Range("I75").Value = 24559
$ws.Range("K75").Value = 24559
$ws.Range("M75").Value = -23623

$ws.Range("H78").Value = 24559
$ws.Range("I78").Value = 24559
$ws.Range("K78").Value = 73677
$ws.Range("M78").Value = -68997

$ws.Range("H107").Value = 1763.2307
$ws.Range("I107").Value = 1844.4
$ws.Range("K107").Value = 5533.200000000001
$ws.Range("M107").Value = -3613.200000000001

$ws.Range("H122").Value = 4022.476
$ws.Range("I122").Value = 3036.7856
$ws.Range("K122").Value = 9110.356800000001
$ws.Range("M122").Value = -6660.356800000001

$ws.Range("H132").Value = 1210.375
$ws.Range("I132").Value = 1183.1666
$ws.Range("J132").Value = 1183.1666
$ws.Range("K132").Value = 3549.4998
$ws.Range("M132").Value = -1019.4998
